$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    3  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    4  = @(0.1190320826869504, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 4.457851494814137)
    5  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    6  = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    7  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    8  = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    9  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    10 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
